# Update cryptos list data (price + volume figures) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.729.42'
$ws.Range("E2").Value = '  +2.98%  '
$ws.Range("D3").Value = '2.562.39'
$ws.Range("E3").Value = '  +4.24%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''501.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.44%  '
$ws.Range("D6").Value = '''152.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.26%  '
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").Value = '''0.576'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.83%  '
$ws.Range("D9").Value = '2.577.69'
$ws.Range("E9").Value = '  +3.72%  '
$ws.Range("E10").Value = '  +3.76%  '
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '''0.339'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").Value = '3.011.17'
$ws.Range("E14").Value = '  +4.73%  '
$ws.Range("D15").Value = '59.948.86'
$ws.Range("E15").Value = '  +3.61%  '
$ws.Range("D16").Value = '''21.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.31%  '
$ws.Range("D17").Value = '''0.0000138'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.93%  '
$ws.Range("D18").Value = '2.573.94'
$ws.Range("E18").Value = '  +4.11%  '
$ws.Range("D19").Value = '''4.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = '''343.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.75%  '
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").Value = '''6.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.57%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").Value = '''59.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.66%  '
$ws.Range("D25").Value = '''0.417'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.28%  '
$ws.Range("D26").Value = '2.700.23'
$ws.Range("E26").Value = '  +6.39%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.80%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '''0.163'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.68%  '
$ws.Range("D29").Value = '0.0₃0842'
$ws.Range("E29").Value = '  +4.16%  '
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").Value = '''155.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("D33").Value = '''19.10'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("E34").Value = '  +0.42%  '
$ws.Range("E35").Value = '  +4.77%  '
$ws.Range("E36").Value = '  +2.94%  '
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("D38").Value = '''0.842'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +22.05%  '
$ws.Range("D39").Value = '''0.836'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("E40").Value = '  +2.71%  '
$ws.Range("D41").Value = '''3.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.68%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '''295.14'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.00%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = '''35.36'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.61%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '''0.616'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.24%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").Value = '''0.0561'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.06%  '
$ws.Range("E46").Value = '  -1.53%  '
$ws.Range("D47").Value = '''0.997'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("D48").Value = '''19.56'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.31%  '
$ws.Range("D49").Value = '''4.88'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.35%  '
$ws.Range("D50").Value = '2.024.99'
$ws.Range("E50").Value = '  +6.38%  '
$ws.Range("D51").Value = '''0.0232'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.10%  '
